$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "vars" (column B) values and new "score" (column C) values for rows 2-51,
# taken from the updated shared-strings table / recomputed scores.
$bVals = @("tejgtotfun_f2prots","tejgfun_f5r08ct05ambpc","_tejgtotfun_f5r08amb","_tejgfun_f1ct05protspc","tejgrb_fonc","_tejgct_r09gstcp","tdvgfun_f1ct06trans","_tejgfun_f5ct06opsegpc","pimgct_r00gstcppc","dfgdevpiagfun_f2ct06agropc","_pimgfun_f5ct06opseg","tejgfun_f5r07ct05cydep","_tejgtotfun_f2opsegpc","_pimgfun_f5r18ct06opseg","mod_110","_dfgdevpiagtotfun_f1pgrco","devppimtotfun_f5r07salud","tdvgge_r00ct06acanfpc","_tejgfun_f5r08ct05opsegpc","_piagtotfun_f5r18amb","dfgpimpiafun_f1ct05prots","pimgfun_f5r18ct06san","_piagfun_f5r07ct05cydep","piagfun_f5r18ct05pgrcopc","_tejgtotfun_f5r18amb","_tdvgge_r00ct05biser","pimgfun_f5r18ct05trans","empinc_2","_dfgdevpiagtotfun_f1trans","piagtotfun_f5r18amb","devppimfun_f5r07ct05agro","_dfgdevpiagfun_f2ct05prsopc","_tejgfun_f5r18ct05pgrco","mod_4","_devppimtotfun_f2amb","dfgpimpiatotfun_f2energpc","pimgfun_f2ct06agropc","_dfgpimpiatotfun_f5r18come","dfgdevpiagge_r08ct05dotra","_dfgpimpiafun_f1ct06transpc","tdvgfun_f1ct06transpc","per_058","_pimgtotfun_f5r18amb","_pimgfun_f5ct05salud","piagfun_f5r07ct06amb","tdvgfun_f2ct06cydep","_devppimfun_f5ct05salud","tdvgfun_f5r18ct05sanpc","tdvgtotfun_f5r07pgrcopc","tdvgfun_f5ct05turi")
$cVals = @(0.03118590079247952,0.02630054578185081,0.02309722825884819,0.01899408735334873,0.01786344312131405,0.0169106237590313,0.01605844311416149,0.01529358886182308,0.01509196776896715,0.01397909969091415,0.01250683888792992,0.01242949441075325,0.01197887491434813,0.01093059033155441,0.01023506559431553,0.01000879146158695,0.009660263545811176,0.009640133008360863,0.008329198695719242,0.008192994631826878,0.007927043363451958,0.007553998846560717,0.00747865904122591,0.007312317378818989,0.006880724802613258,0.006703387945890427,0.006687203887850046,0.006491912994533777,0.006349866278469563,0.006305061746388674,0.006273797247558832,0.006059905979782343,0.006059763487428427,0.005842334590852261,0.00573068531230092,0.005522007588297129,0.005351630039513111,0.005187603179365396,0.005170976277440786,0.005101517308503389,0.004973649512976408,0.004931308794766665,0.004799263551831245,0.004794045351445675,0.004781284369528294,0.004614135250449181,0.00461158249527216,0.004571369849145412,0.004560641013085842,0.00451612425968051)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
}
